$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add/adjust significance stars
$ws.Range("A3").Value = "1.21***"
$ws.Range("B3").Value = "1.25**"
$ws.Range("C3").Value = "1.21**"
$ws.Range("D3").Value = "1.21*"
$ws.Range("E3").Value = "1.21**"

# Row 8: bump star counts on A8, C8, E8 (B8/D8 unchanged)
$ws.Range("A8").Value = "1.46***"
$ws.Range("C8").Value = "1.5***"
$ws.Range("E8").Value = "1.39***"

# Row 18: remove star counts (revert to plain values)
$ws.Range("A18").Value = "1.54"
$ws.Range("B18").Value = "1.57"
$ws.Range("C18").Value = "1.57*"
$ws.Range("D18").Value = "1.5"
$ws.Range("E18").Value = "1.61"

# Row 21: reduce star counts on A21, B21, D21 (C21/E21 unchanged)
$ws.Range("A21").Value = "0.01"
$ws.Range("B21").Value = "18.99**"
$ws.Range("D21").Value = "17.38**"
